# Update rows 3-6 on the "Artfynd" sheet: the four records were
# re-ordered/re-edited (each row's data effectively rotates: old row 6 -> row 3,
# old row 5 -> row 4, old row 3 -> row 5, old row 4 -> row 6), with a few
# field-level edits (added start/end time + public comment text, updated
# reporter/observer names) baked into the new row contents.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (becomes the "Revlummer" / Lycopodium annotinum record)
$ws.Range("A3").Value = 112043158
$ws.Range("B3").Value = 95701
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 221945
$ws.Range("F3").Value = 'Revlummer'
$ws.Range("G3").Value = 'Lycopodium annotinum'
$ws.Range("H3").Value = 'L.'
$ws.Range("P3").Value = 'Stor-Moberg (Stor-Moberg), Dlr'
$ws.Range("Q3").Value = 511628
$ws.Range("R3").Value = 6733623
$ws.Range("S3").Value = 1
$ws.Range("Z3").Value = '10:51'
$ws.Range("AB3").Value = '10:51'
$ws.Range("AC3").Value = 'Finns fläckvis i området'
$ws.Range("AW3").Value = 'Evalena Sköld'
$ws.Range("AX3").Value = 'Evalena Sköld, Åke Sköld'

# Row 4 (becomes the "Blåsippa" / Hepatica nobilis record)
$ws.Range("A4").Value = 112042940
$ws.Range("B4").Value = 98980
$ws.Range("E4").Value = 222498
$ws.Range("F4").Value = 'Blåsippa'
$ws.Range("G4").Value = 'Hepatica nobilis'
$ws.Range("H4").Value = 'Schreb.'
$ws.Range("Q4").Value = 511611
$ws.Range("R4").Value = 6733626
$ws.Range("Z4").Value = '10:33'
$ws.Range("AB4").Value = '10:33'
$ws.Range("AC4").Value = 'Fullt med blåsippsblad på denna sidan bäcken'
$ws.Range("AX4").Value = 'Evalena Sköld, Åke Sköld'

# Row 5 (becomes the "Tretåig hackspett" / Picoides tridactylus record)
$ws.Range("A5").Value = 112042452
$ws.Range("B5").Value = 56430
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = 'Tretåig hackspett'
$ws.Range("G5").Value = 'Picoides tridactylus'
$ws.Range("H5").Value = '(Linnaeus, 1758)'
$ws.Range("P5").Value = 'Stor Mpmerg, Kilen-Stor, Moberg, Leksand, Dlr'
$ws.Range("Q5").Value = 511614
$ws.Range("R5").Value = 6733640
$ws.Range("S5").Value = 25
$ws.Range("Z5").Value = $null
$ws.Range("AB5").Value = $null
$ws.Range("AC5").Value = $null
$ws.Range("AW5").Value = 'Åke Sköld'
$ws.Range("AX5").Value = 'Åke Sköld'

# Row 6 (becomes the "Svavelriska" / Lactarius scrobiculatus record)
$ws.Range("A6").Value = 112043031
$ws.Range("B6").Value = 90480
$ws.Range("E6").Value = 4769
$ws.Range("F6").Value = 'Svavelriska'
$ws.Range("G6").Value = 'Lactarius scrobiculatus'
$ws.Range("H6").Value = '(Scop.:Fr.) Fr.'
$ws.Range("Q6").Value = 511625
$ws.Range("R6").Value = 6733616
$ws.Range("Z6").Value = '10:42'
$ws.Range("AB6").Value = '10:42'
$ws.Range("AC6").Value = $null
$ws.Range("AX6").Value = 'Evalena Sköld'
